$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172834515571594
$ws.Range("B1").Value = 2.154973268508911
$ws.Range("C1").Value = 3.153638362884521
$ws.Range("D1").Value = 3.685949325561523
$ws.Range("E1").Value = 1.371334910392761
